$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45178 -> 45179) for every data row (rows 2 through 97).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45178) {
        $cell.Value = 45179
    }
}
